{"js": "// Update the 20x5 arithmetic-answers table: replace each cell's equation\n// text with the new equation per the commit diff, preserving all\n// paragraph/run formatting (fonts, size, alignment) already on the cells.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\"25+45=70\", \"26+40=66\", \"46+45=91\", \"24+70=94\", \"5+37=42\"],\n  [\"42+21=63\", \"97-59=38\", \"8+51=59\", \"38+1=39\", \"21+52=73\"],\n  [\"46-29=17\", \"30-12=18\", \"41+15=56\", \"6+22=28\", \"12+73=85\"],\n  [\"19+79=98\", \"16-15=1\", \"21-9=12\", \"44+13=57\", \"46+37=83\"],\n  [\"32-21=11\", \"42+55=97\", \"6+17=23\", \"11-7=4\", \"71-20=51\"],\n  [\"9+84=93\", \"32-29=3\", \"33-20=13\", \"98-26=72\", \"27+33=60\"],\n  [\"69+10=79\", \"96-75=21\", \"78+4=82\", \"97-36=61\", \"69+24=93\"],\n  [\"11+20=31\", \"64-13=51\", \"48+51=99\", \"79-48=31\", \"96-6=90\"],\n  [\"44+20=64\", \"30+36=66\", \"90-76=14\", \"42-8=34\", \"19+18=37\"],\n  [\"36+9=45\", \"57+41=98\", \"28-28=0\", \"22+7=29\", \"73-6=67\"],\n  [\"30+57=87\", \"64+35=99\", \"55-43=12\", \"67-2=65\", \"87-72=15\"],\n  [\"74+12=86\", \"35+40=75\", \"68-58=10\", \"93-78=15\", \"40+22=62\"],\n  [\"3+34=37\", \"7+41=48\", \"40-39=1\", \"16-1=15\", \"71-49=22\"],\n  [\"97-31=66\", \"42+41=83\", \"9+43=52\", \"37+21=58\", \"71-52=19\"],\n  [\"53-14=39\", \"7+4=11\", \"72-15=57\", \"60-24=36\", \"8+43=51\"],\n  [\"44+9=53\", \"25-24=1\", \"7+20=27\", \"77-72=5\", \"44+23=67\"],\n  [\"38+56=94\", \"34+32=66\", \"19+27=46\", \"85-34=51\", \"8+89=97\"],\n  [\"13+55=68\", \"54-47=7\", \"41-17=24\", \"30+30=60\", \"61+10=71\"],\n  [\"95-38=57\", \"56-2=54\", \"92-31=61\", \"51+5=56\", \"5+40=45\"],\n  [\"95-69=26\", \"48+0=48\", \"82-51=31\", \"27+39=66\", \"96-87=9\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the 20x5 arithmetic-answers table: replace each cell's equation\n# text with the new equation per the commit diff. Setting Cell.Range.Text\n# replaces only the cell's text content (Word keeps the trailing cell\n# mark / paragraph mark automatically) so the existing run formatting\n# (TimeNewRoman font, size 30, left alignment) on each cell is preserved.\n\n$newValues = @(\n    @(\"25+45=70\", \"26+40=66\", \"46+45=91\", \"24+70=94\", \"5+37=42\"),\n    @(\"42+21=63\", \"97-59=38\", \"8+51=59\", \"38+1=39\", \"21+52=73\"),\n    @(\"46-29=17\", \"30-12=18\", \"41+15=56\", \"6+22=28\", \"12+73=85\"),\n    @(\"19+79=98\", \"16-15=1\", \"21-9=12\", \"44+13=57\", \"46+37=83\"),\n    @(\"32-21=11\", \"42+55=97\", \"6+17=23\", \"11-7=4\", \"71-20=51\"),\n    @(\"9+84=93\", \"32-29=3\", \"33-20=13\", \"98-26=72\", \"27+33=60\"),\n    @(\"69+10=79\", \"96-75=21\", \"78+4=82\", \"97-36=61\", \"69+24=93\"),\n    @(\"11+20=31\", \"64-13=51\", \"48+51=99\", \"79-48=31\", \"96-6=90\"),\n    @(\"44+20=64\", \"30+36=66\", \"90-76=14\", \"42-8=34\", \"19+18=37\"),\n    @(\"36+9=45\", \"57+41=98\", \"28-28=0\", \"22+7=29\", \"73-6=67\"),\n    @(\"30+57=87\", \"64+35=99\", \"55-43=12\", \"67-2=65\", \"87-72=15\"),\n    @(\"74+12=86\", \"35+40=75\", \"68-58=10\", \"93-78=15\", \"40+22=62\"),\n    @(\"3+34=37\", \"7+41=48\", \"40-39=1\", \"16-1=15\", \"71-49=22\"),\n    @(\"97-31=66\", \"42+41=83\", \"9+43=52\", \"37+21=58\", \"71-52=19\"),\n    @(\"53-14=39\", \"7+4=11\", \"72-15=57\", \"60-24=36\", \"8+43=51\"),\n    @(\"44+9=53\", \"25-24=1\", \"7+20=27\", \"77-72=5\", \"44+23=67\"),\n    @(\"38+56=94\", \"34+32=66\", \"19+27=46\", \"85-34=51\", \"8+89=97\"),\n    @(\"13+55=68\", \"54-47=7\", \"41-17=24\", \"30+30=60\", \"61+10=71\"),\n    @(\"95-38=57\", \"56-2=54\", \"92-31=61\", \"51+5=56\", \"5+40=45\"),\n    @(\"95-69=26\", \"48+0=48\", \"82-51=31\", \"27+39=66\", \"96-87=9\"),\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowVals = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowVals.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
